# DataEngine.xlsx edit: ignore string upper/lower case
# - Clears stale PASS/FAIL "Results" markers from most TestSteps rows
# - Updates the first TestSteps row's ActionKeyword/TestData/Results values
# - Updates a leftover "8" TestData value to "username"
# - Moves the active sheet / selection from TestCases to TestSteps

$wb  = $excel.ActiveWorkbook
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestSteps = $wb.Worksheets.Item("TestSteps")

# --- Cell value edits on the TestSteps sheet ---------------------------

# Row 2: TS_1 / open browser step
$wsTestSteps.Range("E2").Value = "eat"
$wsTestSteps.Range("F2").Value = "Browser"
$wsTestSteps.Range("G2").Value = "FAIL"

# Rows 3-6 and 8: drop the stale "PASS" result marker entirely
$wsTestSteps.Range("G3").ClearContents()
$wsTestSteps.Range("G4").ClearContents()
$wsTestSteps.Range("G5").ClearContents()
$wsTestSteps.Range("G6").ClearContents()
$wsTestSteps.Range("G8").ClearContents()

# Row 7: TestData changes from "8" to "username"; drop its "FAIL" marker
$wsTestSteps.Range("F7").Value = "username"
$wsTestSteps.Range("G7").ClearContents()

# --- Sheet view / active-tab changes ------------------------------------

# Active tab moves from TestCases (index 0) to TestSteps (index 1)
[void]$wsTestSteps.Activate()

# TestCases loses its "selected tab" view state; selection moves to B13
[void]$wsTestCases.Range("B13").Select()

# TestSteps becomes the selected tab; selection moves to D20
[void]$wsTestSteps.Range("D20").Select()
